# Novi.xlsx update: append daily-data rows through 2022-01-05 (commit:
# "aggiornamento fino a 6 gennaio 2022"), extending the sheet from row 464
# (date serial 44538) through row 491 (date serial 44566).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, date (serial), col B, col C, col D
$data = @(
  @(465, 44539, 4, 14, 141.9590346785642),
  @(466, 44540, 4, 16, 162.2388967755019),
  @(467, 44541, 0, 13, 131.8191036300953),
  @(468, 44542, 1, 12, 121.6791725816264),
  @(469, 44543, 5, 17, 172.3788278239708),
  @(470, 44544, 7, 22, 223.0784830663152),
  @(471, 44545, 0, 21, 212.9385520178463),
  @(472, 44546, 0, 17, 172.3788278239708),
  @(473, 44547, 9, 22, 223.0784830663152),
  @(474, 44548, 1, 23, 233.218414114784),
  @(475, 44550, 4, 26, 263.6382072601907),
  @(476, 44551, 7, 28, 283.9180693571284),
  @(477, 44552, 0, 21, 212.9385520178463),
  @(478, 44553, 5, 26, 263.6382072601907),
  @(479, 44554, 1, 27, 273.7781383086595),
  @(480, 44555, 8, 26, 263.6382072601907),
  @(481, 44556, 2, 27, 273.7781383086595),
  @(482, 44557, 6, 29, 294.0580004055972),
  @(483, 44558, 2, 24, 243.3583451632529),
  @(484, 44559, 17, 41, 415.7371729872237),
  @(485, 44560, 1, 37, 375.1774487933482),
  @(486, 44561, 5, 41, 415.7371729872237),
  @(487, 44562, 6, 39, 395.4573108902859),
  @(488, 44563, 21, 58, 588.1160008111945),
  @(489, 44564, 14, 66, 669.2354491989455),
  @(490, 44565, 5, 69, 699.655242344352),
  @(491, 44566, 8, 60, 608.3958629081322)
)

# The last existing row (464) carries the date-column cell style (centered,
# bordered, "YYYY-MM-DD HH:MM:SS" number format). Propagate that formatting
# to column A of every new row by copying formats only (keeps the same
# style index instead of minting a new one).
$ws.Range("A464").Copy() | Out-Null

foreach ($row in $data) {
    $r = [int]$row[0]

    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($r, 1).Value2 = $row[1]
    $ws.Cells.Item($r, 2).Value2 = $row[2]
    $ws.Cells.Item($r, 3).Value2 = $row[3]
    $ws.Cells.Item($r, 4).Value2 = $row[4]
}

$excel.CutCopyMode = $false
